# Daily attendance processing - 2026-01-20 11:39:23
#
# The "Recorded By" column (G) lists the people/agents who recorded each
# attendance session, separated by commas. For every row where it reads
# "dnasr281@gmail.com, System" the order of the two entries is flipped to
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Restrict the search to column G so we never touch any other column,
# and use Find/FindNext so untouched (empty) cells are never read.
$searchRange = $ws.Range("G1:G$lastRow")

$firstMatch = $searchRange.Find($oldValue)
if ($firstMatch -ne $null) {
    $firstAddress = $firstMatch.Address()
    $current = $firstMatch
    $guard = 0
    do {
        $current.Value = $newValue
        $current = $searchRange.FindNext($current)
        $guard++
    } while ($current -ne $null -and $current.Address() -ne $firstAddress -and $guard -lt 1000)
}
